# Update LR-pair TPM-derived metrics with new TPM-based calculations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2115286666666667
$ws.Range("H2").Value = 0.634586
$ws.Range("I2").Value = 0.08153347995807345
$ws.Range("J2").Value = 0.08153347995807345
$ws.Range("M2").Value = 0.140567
$ws.Range("N2").Value = 0.421701
$ws.Range("O2").Value = 0.07810038533383065
$ws.Range("P2").Value = 0.07810038533383065
$ws.Range("Q2").Value = 0.02973395008733333
$ws.Range("R2").Value = 0.267605550786
$ws.Range("S2").Value = 0.006367796202333695
$ws.Range("T2").Value = 0.006367796202333695

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2115286666666667
$ws.Range("H3").Value = 0.634586
$ws.Range("I3").Value = 0.08153347995807345
$ws.Range("J3").Value = 0.08153347995807345
$ws.Range("O3").Value = 0.02984383293631935
$ws.Range("P3").Value = 0.02984383293631935
$ws.Range("Q3").Value = 0.01136198029177778
$ws.Range("R3").Value = 0.102257822626
$ws.Range("S3").Value = 0.002433271554585486
$ws.Range("T3").Value = 0.002433271554585486

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2115286666666667
$ws.Range("H4").Value = 0.634586
$ws.Range("I4").Value = 0.08153347995807345
$ws.Range("J4").Value = 0.08153347995807345
$ws.Range("O4").Value = 0.8920557817298499
$ws.Range("P4").Value = 0.8920557817298499
$ws.Range("Q4").Value = 0.3396185815946667
$ws.Range("R4").Value = 3.056567234352
$ws.Range("S4").Value = 0.07273241220115426
$ws.Range("T4").Value = 0.07273241220115426

# Row 5
$ws.Range("I5").Value = 0.4735790235655714
$ws.Range("J5").Value = 0.4735790235655714
$ws.Range("M5").Value = 0.140567
$ws.Range("N5").Value = 0.421701
$ws.Range("O5").Value = 0.07810038533383065
$ws.Range("P5").Value = 0.07810038533383065
$ws.Range("Q5").Value = 0.172706660581
$ws.Range("R5").Value = 1.554359945229
$ws.Range("S5").Value = 0.03698670422649039
$ws.Range("T5").Value = 0.03698670422649039

# Row 6
$ws.Range("I6").Value = 0.4735790235655714
$ws.Range("J6").Value = 0.4735790235655714
$ws.Range("O6").Value = 0.02984383293631935
$ws.Range("P6").Value = 0.02984383293631935
$ws.Range("S6").Value = 0.01413341326143616
$ws.Range("T6").Value = 0.01413341326143616

# Row 7
$ws.Range("I7").Value = 0.4735790235655714
$ws.Range("J7").Value = 0.4735790235655714
$ws.Range("O7").Value = 0.8920557817298499
$ws.Range("P7").Value = 0.8920557817298499
$ws.Range("S7").Value = 0.4224589060776448
$ws.Range("T7").Value = 0.4224589060776448

# Row 8
$ws.Range("I8").Value = 0.4448874964763552
$ws.Range("J8").Value = 0.4448874964763552
$ws.Range("M8").Value = 0.140567
$ws.Range("N8").Value = 0.421701
$ws.Range("O8").Value = 0.07810038533383065
$ws.Range("P8").Value = 0.07810038533383065
$ws.Range("Q8").Value = 0.1622433216576667
$ws.Range("R8").Value = 1.460189894919
$ws.Range("S8").Value = 0.03474588490500657
$ws.Range("T8").Value = 0.03474588490500657

# Row 9
$ws.Range("I9").Value = 0.4448874964763552
$ws.Range("J9").Value = 0.4448874964763552
$ws.Range("O9").Value = 0.02984383293631935
$ws.Range("P9").Value = 0.02984383293631935
$ws.Range("S9").Value = 0.01327714812029771
$ws.Range("T9").Value = 0.01327714812029771

# Row 10
$ws.Range("I10").Value = 0.4448874964763552
$ws.Range("J10").Value = 0.4448874964763552
$ws.Range("O10").Value = 0.8920557817298499
$ws.Range("P10").Value = 0.8920557817298499
$ws.Range("S10").Value = 0.3968644634510509
$ws.Range("T10").Value = 0.3968644634510509
